$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) column C was refreshed for every existing data
# row (2-350): 45202 -> 45203 (2023-10-03 -> 2023-10-04).
$ws.Range("C2:C350").Value = 45203

# Row 350 picks up an explicit row height (15pt) in the new file.
$ws.Rows.Item(350).RowHeight = 15

# A brand-new entry (row 351) was appended for case "A 47282-2023".
# Seed formatting by copying from row 350 (skipping the unused column F),
# then overwrite the values.
$ws.Range("A350:E350").Copy($ws.Range("A351:E351"))
$ws.Range("G350:R350").Copy($ws.Range("G351:R351"))

$ws.Range("A351").Value = "A 47282-2023"
$ws.Range("B351").Value = 45202
$ws.Range("C351").Value = 45203
$ws.Range("D351").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E351").Value = "KÖPING"
$ws.Range("G351").Value = 5
$ws.Range("H351").Value = 0
$ws.Range("I351").Value = 0
$ws.Range("J351").Value = 0
$ws.Range("K351").Value = 0
$ws.Range("L351").Value = 0
$ws.Range("M351").Value = 0
$ws.Range("N351").Value = 0
$ws.Range("O351").Value = 0
$ws.Range("P351").Value = 0
$ws.Range("Q351").Value = 0
$ws.Range("R351").Value = ""
